$d = $word.ActiveDocument

# 1. Remove the "Meta description" paragraph (the second paragraph of the
#    document) entirely, including its paragraph mark.
$metaPara = $d.Paragraphs(2)
[void]$metaPara.Range.Delete()

# 2. Replace the text of the final paragraph (the italic "Create a
#    cartoon-style feature image..." paragraph) with the meta-description
#    copy, keeping its existing italic run formatting untouched.
[void]$d.Content.Find.Execute(
    "Create a cartoon-style feature image for Atlantis Megaways that features a happy Maya warrior wearing glasses. The warrior should be positioned underwater among ruins of the lost city of Atlantis with sea creatures swimming around in the background. The image should incorporate the game's logo and feature vibrant colors that capture the adventurous and mysterious theme of the game. The image should also clearly convey the idea of winning cash prizes with a bubbly, celebratory vibe.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Read our Atlantis Megaways slot game review and play for free. Features, gameplay mechanics, betting range, and jackpot prizes explained.",
    2)

# 3. Insert a new bold paragraph just before that (now updated) paragraph,
#    holding the page title copy that used to live in the Meta description
#    paragraph near the top of the document.
$lastIndex = $d.Paragraphs.Count
$targetPara = $d.Paragraphs($lastIndex)
[void]$targetPara.Range.InsertParagraphBefore()

$newParaXml = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Atlantis Megaways Slot Game for Free | Review</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$insertedPara = $d.Paragraphs($lastIndex)
[void]$insertedPara.Range.InsertXML($newParaXml)
